# Auto-generated Excel COM-interop script to apply market-data value updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets (FFXIV Leve profit data).
$wb = $excel.ActiveWorkbook

# ---- Worksheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 78.22221999999999
$ws.Cells.Item(6, 9).Value = 78.22221999999999
$ws.Cells.Item(6, 11).Value = 234.66666
$ws.Cells.Item(6, 13).Value = -122.66666
$ws.Cells.Item(7, 8).Value = 772.75
$ws.Cells.Item(7, 10).Value = 550.5
$ws.Cells.Item(7, 12).Value = 550.5
$ws.Cells.Item(7, 14).Value = -774.5
$ws.Cells.Item(8, 8).Value = 122.888885
$ws.Cells.Item(8, 9).Value = 122.888885
$ws.Cells.Item(8, 11).Value = 368.666655
$ws.Cells.Item(8, 13).Value = -229.666655
$ws.Cells.Item(9, 8).Value = 136.55556
$ws.Cells.Item(9, 9).Value = 203.25
$ws.Cells.Item(9, 10).Value = 83.2
$ws.Cells.Item(9, 11).Value = 203.25
$ws.Cells.Item(9, 12).Value = 83.2
$ws.Cells.Item(9, 13).Value = -34.25
$ws.Cells.Item(9, 14).Value = -421.2
$ws.Cells.Item(10, 8).Value = 1500
$ws.Cells.Item(10, 10).Value = 2000
$ws.Cells.Item(10, 12).Value = 2000
$ws.Cells.Item(10, 14).Value = -2586
$ws.Cells.Item(13, 8).Value = 1167
$ws.Cells.Item(13, 10).Value = 2500
$ws.Cells.Item(13, 12).Value = 2500
$ws.Cells.Item(13, 14).Value = -2838
$ws.Cells.Item(14, 8).Value = 772.75
$ws.Cells.Item(14, 10).Value = 550.5
$ws.Cells.Item(14, 12).Value = 550.5
$ws.Cells.Item(14, 14).Value = -932.5
$ws.Cells.Item(16, 8).Value = 6716.6665
$ws.Cells.Item(16, 9).Value = 50
$ws.Cells.Item(16, 10).Value = 10050
$ws.Cells.Item(16, 11).Value = 50
$ws.Cells.Item(16, 12).Value = 10050
$ws.Cells.Item(16, 13).Value = 180
$ws.Cells.Item(16, 14).Value = -10510
$ws.Cells.Item(17, 8).Value = 914.03705
$ws.Cells.Item(17, 10).Value = 912.0961
$ws.Cells.Item(17, 12).Value = 2736.2883
$ws.Cells.Item(17, 14).Value = -3072.2883
$ws.Cells.Item(76, 8).Value = 3857.1428
$ws.Cells.Item(79, 8).Value = 3857.1428
$ws.Cells.Item(80, 8).Value = 6580109.5
$ws.Cells.Item(80, 9).Value = 9616392
$ws.Cells.Item(80, 10).Value = 1498
$ws.Cells.Item(80, 11).Value = 28849176
$ws.Cells.Item(80, 12).Value = 4494
$ws.Cells.Item(80, 13).Value = -28848178
$ws.Cells.Item(80, 14).Value = -6490
$ws.Cells.Item(83, 8).Value = 6580109.5
$ws.Cells.Item(83, 9).Value = 9616392
$ws.Cells.Item(83, 10).Value = 1498
$ws.Cells.Item(83, 11).Value = 86547528
$ws.Cells.Item(83, 12).Value = 13482
$ws.Cells.Item(83, 13).Value = -86542536
$ws.Cells.Item(83, 14).Value = -23466
$ws.Cells.Item(96, 8).Value = 2832.7
$ws.Cells.Item(96, 10).Value = 5114.25
$ws.Cells.Item(96, 12).Value = 15342.75
$ws.Cells.Item(96, 14).Value = -18088.75
$ws.Cells.Item(106, 8).Value = 4493.1665
$ws.Cells.Item(106, 9).Value = 4425.4443
$ws.Cells.Item(106, 10).Value = 4696.3335
$ws.Cells.Item(106, 11).Value = 4425.4443
$ws.Cells.Item(106, 12).Value = 4696.3335
$ws.Cells.Item(106, 13).Value = -3794.4443
$ws.Cells.Item(106, 14).Value = -5958.3335
$ws.Cells.Item(111, 8).Value = 4273.091
$ws.Cells.Item(111, 9).Value = 3529.8333
$ws.Cells.Item(111, 11).Value = 10589.4999
$ws.Cells.Item(111, 13).Value = -7522.499899999999
$ws.Cells.Item(113, 8).Value = 3549.5
$ws.Cells.Item(113, 9).Value = 3310.4443
$ws.Cells.Item(113, 11).Value = 3310.4443
$ws.Cells.Item(113, 13).Value = -56.44430000000011

# ---- Worksheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2814.8572
$ws.Cells.Item(2, 9).Value = 1620.6111
$ws.Cells.Item(2, 11).Value = 1620.6111
$ws.Cells.Item(2, 13).Value = -1507.6111
$ws.Cells.Item(116, 8).Value = 2814.8572
$ws.Cells.Item(116, 9).Value = 1620.6111
$ws.Cells.Item(116, 11).Value = 1620.6111
$ws.Cells.Item(116, 13).Value = 673.3888999999999
$ws.Cells.Item(132, 8).Value = 767362.0600000001
$ws.Cells.Item(132, 9).Value = 913766.5600000001
$ws.Cells.Item(132, 11).Value = 2741299.68
$ws.Cells.Item(132, 13).Value = -2738769.68

# ---- Worksheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2814.8572
$ws.Cells.Item(3, 9).Value = 1620.6111
$ws.Cells.Item(3, 11).Value = 1620.6111
$ws.Cells.Item(3, 13).Value = -1506.6111

# ---- Worksheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 27293700
$ws.Cells.Item(132, 10).Value = 13425497
$ws.Cells.Item(132, 12).Value = 40276491
$ws.Cells.Item(132, 14).Value = -40281551
$ws.Cells.Item(134, 8).Value = 3276938.8
$ws.Cells.Item(134, 9).Value = 18936.37
$ws.Cells.Item(134, 11).Value = 56809.11
$ws.Cells.Item(134, 13).Value = -54274.11

# ---- Worksheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(6, 8).Value = 323
$ws.Cells.Item(6, 9).Value = 247.14285
$ws.Cells.Item(6, 11).Value = 741.4285500000001
$ws.Cells.Item(6, 13).Value = -628.4285500000001
$ws.Cells.Item(113, 8).Value = 3746.8
$ws.Cells.Item(113, 10).Value = 5938.3335
$ws.Cells.Item(113, 12).Value = 17815.0005
$ws.Cells.Item(113, 14).Value = -22155.0005
$ws.Cells.Item(132, 8).Value = 2028.5714
$ws.Cells.Item(132, 9).Value = 1475.25
$ws.Cells.Item(132, 11).Value = 13277.25
$ws.Cells.Item(132, 13).Value = -10747.25

# ---- Worksheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(6, 8).Value = 925
$ws.Cells.Item(6, 9).Value = 1000
$ws.Cells.Item(6, 10).Value = 850
$ws.Cells.Item(6, 11).Value = 1000
$ws.Cells.Item(6, 12).Value = 850
$ws.Cells.Item(6, 13).Value = -887
$ws.Cells.Item(6, 14).Value = -1076
$ws.Cells.Item(16, 8).Value = 925
$ws.Cells.Item(16, 9).Value = 1000
$ws.Cells.Item(16, 10).Value = 850
$ws.Cells.Item(16, 11).Value = 1000
$ws.Cells.Item(16, 12).Value = 850
$ws.Cells.Item(16, 13).Value = -750
$ws.Cells.Item(16, 14).Value = -1350
$ws.Cells.Item(31, 8).Value = 1058.8
$ws.Cells.Item(31, 9).Value = 1058.8
$ws.Cells.Item(31, 11).Value = 1058.8
$ws.Cells.Item(31, 13).Value = -766.8
$ws.Cells.Item(37, 8).Value = 1058.8
$ws.Cells.Item(37, 9).Value = 1058.8
$ws.Cells.Item(37, 11).Value = 1058.8
$ws.Cells.Item(37, 13).Value = -781.8
$ws.Cells.Item(43, 8).Value = 3396.077
$ws.Cells.Item(43, 9).Value = 1418
$ws.Cells.Item(43, 10).Value = 14275.5
$ws.Cells.Item(43, 11).Value = 1418
$ws.Cells.Item(43, 12).Value = 14275.5
$ws.Cells.Item(43, 13).Value = -1267
$ws.Cells.Item(43, 14).Value = -14577.5
$ws.Cells.Item(80, 8).Value = 2324.7273
$ws.Cells.Item(80, 9).Value = 2419.4443
$ws.Cells.Item(80, 10).Value = 2211.0667
$ws.Cells.Item(80, 11).Value = 2419.4443
$ws.Cells.Item(80, 12).Value = 2211.0667
$ws.Cells.Item(80, 13).Value = -1421.4443
$ws.Cells.Item(80, 14).Value = -4207.066699999999
$ws.Cells.Item(83, 8).Value = 2324.7273
$ws.Cells.Item(83, 9).Value = 2419.4443
$ws.Cells.Item(83, 10).Value = 2211.0667
$ws.Cells.Item(83, 11).Value = 12097.2215
$ws.Cells.Item(83, 12).Value = 11055.3335
$ws.Cells.Item(83, 13).Value = -7105.2215
$ws.Cells.Item(83, 14).Value = -21039.3335
$ws.Cells.Item(97, 8).Value = 1463.2593
$ws.Cells.Item(97, 9).Value = 1504.75
$ws.Cells.Item(97, 11).Value = 1504.75
$ws.Cells.Item(97, 13).Value = -1008.75
$ws.Cells.Item(122, 8).Value = 31987.621
$ws.Cells.Item(122, 9).Value = 46125.78
$ws.Cells.Item(122, 11).Value = 138377.34
$ws.Cells.Item(122, 13).Value = -135927.34

# ---- Worksheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 220.45454
$ws.Cells.Item(9, 9).Value = 142.6
$ws.Cells.Item(9, 10).Value = 999
$ws.Cells.Item(9, 11).Value = 142.6
$ws.Cells.Item(9, 12).Value = 999
$ws.Cells.Item(9, 13).Value = 81.40000000000001
$ws.Cells.Item(9, 14).Value = -1447
$ws.Cells.Item(39, 8).Value = 17499
$ws.Cells.Item(39, 9).Value = 0
$ws.Cells.Item(39, 10).Value = 17499
$ws.Cells.Item(39, 11).Value = 0
$ws.Cells.Item(39, 12).Value = 17499
$ws.Cells.Item(39, 13).ClearContents()
$ws.Cells.Item(39, 14).Value = -18419
$ws.Cells.Item(45, 8).Value = 0
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 13).ClearContents()
$ws.Cells.Item(48, 8).Value = 34999
$ws.Cells.Item(48, 10).Value = 34999
$ws.Cells.Item(48, 12).Value = 34999
$ws.Cells.Item(48, 14).Value = -36321
$ws.Cells.Item(61, 8).Value = 2059.4
$ws.Cells.Item(61, 9).Value = 1949.25
$ws.Cells.Item(61, 11).Value = 1949.25
$ws.Cells.Item(61, 13).Value = -1747.25
$ws.Cells.Item(82, 8).Value = 879.8946999999999
$ws.Cells.Item(82, 9).Value = 617.2
$ws.Cells.Item(82, 11).Value = 617.2
$ws.Cells.Item(82, 13).Value = -256.2
$ws.Cells.Item(85, 8).Value = 879.8946999999999
$ws.Cells.Item(85, 9).Value = 617.2
$ws.Cells.Item(85, 11).Value = 617.2
$ws.Cells.Item(85, 13).Value = 630.8
$ws.Cells.Item(113, 8).Value = 2059.4
$ws.Cells.Item(113, 9).Value = 1949.25
$ws.Cells.Item(113, 11).Value = 1949.25
$ws.Cells.Item(113, 13).Value = 220.75

# ---- Worksheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(75, 8).Value = 32104
$ws.Cells.Item(75, 10).Value = 25130
$ws.Cells.Item(75, 12).Value = 25130
$ws.Cells.Item(75, 14).Value = -27002
$ws.Cells.Item(78, 8).Value = 32104
$ws.Cells.Item(78, 10).Value = 25130
$ws.Cells.Item(78, 12).Value = 75390
$ws.Cells.Item(78, 14).Value = -84750
$ws.Cells.Item(113, 8).Value = 1474.5358
$ws.Cells.Item(113, 9).Value = 293.66666
$ws.Cells.Item(113, 10).Value = 2837.077
$ws.Cells.Item(113, 11).Value = 880.9999799999999
$ws.Cells.Item(113, 12).Value = 8511.231
$ws.Cells.Item(113, 13).Value = 1289.00002
$ws.Cells.Item(113, 14).Value = -12851.231
$ws.Cells.Item(122, 8).Value = 2453.7844
$ws.Cells.Item(122, 9).Value = 2138.348
$ws.Cells.Item(122, 11).Value = 6415.044
$ws.Cells.Item(122, 13).Value = -3965.044
$ws.Cells.Item(132, 8).Value = 11148409
$ws.Cells.Item(132, 9).Value = 14862374
$ws.Cells.Item(132, 11).Value = 44587122
$ws.Cells.Item(132, 13).Value = -44584592
